$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.901.40"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "3.365.16"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'405.66"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").Value = "'133.39"
$ws.Range("E6").Value = "  +9.25%  "
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.670"
$ws.Range("E9").Value = "  +4.26%  "
$ws.Range("E10").Value = "  +7.01%  "
$ws.Range("D11").Value = "'42.27"
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "3.885.91"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").Value = "'8.31"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "3.349.86"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "60.973.65"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "'11.06"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("E20").Value = "  +6.80%  "
$ws.Range("D21").Value = "'3.21"
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").Value = "'83.60"
$ws.Range("E22").Value = "  +9.77%  "
$ws.Range("D23").Value = "'12.70"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Value = "'304.20"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "'3.12"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "'4.77"
$ws.Range("E26").Value = "  +11.92%  "
$ws.Range("D27").Value = "'8.32"
$ws.Range("E27").Value = "  +8.35%  "
$ws.Range("D28").Value = "'29.34"
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("D29").Value = "'7.55"
$ws.Range("E29").Value = "  -7.71%  "
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'11.29"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").Value = "'41.12"
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("D35").Value = "'2.47"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("D37").Value = "'51.71"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  -3.83%  "
$ws.Range("D40").Value = "'2.89"
$ws.Range("E40").Value = "  -4.29%  "
$ws.Range("D41").Value = "'1.98"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "'137.22"
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "'4.00"
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("D46").Value = "'16.62"
$ws.Range("E46").Value = "  -4.38%  "
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("D48").Value = "'21.22"
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("D49").Value = "2.121.00"
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").Value = "'2.30"
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("D51").Value = "'1.88"
$ws.Range("E51").Value = "  -0.09%  "
